# Insert a new data row at row 242 (pushes existing rows 242-261 down to 243-262)
# and populate it with the new observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(242).Insert()

$ws.Cells.Item(242, 1).Value = 5
$ws.Cells.Item(242, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(242, 3).Value = "Maule"
$ws.Cells.Item(242, 4).Value = 44746
$ws.Cells.Item(242, 5).Value = 7
$ws.Cells.Item(242, 6).Value = 100112045
$ws.Cells.Item(242, 7).Value = "Zapallo"
$ws.Cells.Item(242, 8).Value = "Camote"
$ws.Cells.Item(242, 9).Value = "1a (guarda)"
$ws.Cells.Item(242, 10).Value = 900
$ws.Cells.Item(242, 11).Value = 400
$ws.Cells.Item(242, 12).Value = 400
$ws.Cells.Item(242, 13).Value = 400
$ws.Cells.Item(242, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(242, 15).Value = "Región del Maule"
$ws.Cells.Item(242, 16).Value = 400
$ws.Cells.Item(242, 17).Value = 1
$ws.Cells.Item(242, 18).Value = "Hortaliza"
